$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = -0.039877723698143519
$ws.Range("B1").Value = 0.039877722998931976

# Row 2
$ws.Range("A2").Value = 0.022169434153850157
$ws.Range("B2").Value = -0.022169434891748163

# Row 3
$ws.Range("A3").Value = -0.028435964436284107
$ws.Range("B3").Value = 0.028435963728436469

# Row 4 (new)
$ws.Range("A4").Value = -0.044802605598396147
$ws.Range("B4").Value = 0.04480260487590261

# Row 5 (new)
$ws.Range("A5").Value = 0.060155040415389625
$ws.Range("B5").Value = -0.060155041140464667

# Column A width now matches column B's width (14.42578125 characters).
# The COM layer quantizes ColumnWidth to a coarse grid, so this is the
# closest achievable setting that lands on the target width (14.5).
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
